$wb = $excel.ActiveWorkbook

# Values updated by the scheduled Kraken_Profits profit-recalculation runner.
# Each worksheet corresponds to a crafting job (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR);
# cells H:N hold market-price / profit figures that the runner refreshes per leve row.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 261.05554
$ws.Range("I9").Value = 184.4
$ws.Range("J9").Value = 644.3333
$ws.Range("K9").Value = 184.4
$ws.Range("L9").Value = 644.3333
$ws.Range("M9").Value = -15.40000000000001
$ws.Range("N9").Value = -982.3333
$ws.Range("H18").Value = 11641.143
$ws.Range("J18").Value = 1197.8
$ws.Range("L18").Value = 1197.8
$ws.Range("N18").Value = -1765.8
$ws.Range("H63").Value = 19000
$ws.Range("I63").Value = 19000
$ws.Range("K63").Value = 19000
$ws.Range("M63").Value = -18376
$ws.Range("H66").Value = 19000
$ws.Range("I66").Value = 19000
$ws.Range("K66").Value = 57000
$ws.Range("M66").Value = -53880
$ws.Range("H68").Value = 29000
$ws.Range("J68").Value = 39000
$ws.Range("L68").Value = 39000
$ws.Range("N68").Value = -40498
$ws.Range("H71").Value = 29000
$ws.Range("J71").Value = 39000
$ws.Range("L71").Value = 117000
$ws.Range("N71").Value = -124488
$ws.Range("H92").Value = 1215.5
$ws.Range("J92").Value = 2747.5
$ws.Range("L92").Value = 2747.5
$ws.Range("N92").Value = -5243.5
$ws.Range("H106").Value = 200
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("H115").Value = 290
$ws.Range("I115").Value = 290
$ws.Range("K115").Value = 870
$ws.Range("M115").Value = 697
$ws.Range("H125").Value = 1000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 9000
$ws.Range("N125").Value = -13920
$ws.Range("H132").Value = 3858.647
$ws.Range("I132").Value = 3207.1538
$ws.Range("J132").Value = 5976
$ws.Range("K132").Value = 9621.4614
$ws.Range("L132").Value = 17928
$ws.Range("M132").Value = -7091.4614
$ws.Range("N132").Value = -22988
$ws.Range("H141").Value = 38974.25
$ws.Range("I141").Value = 7949
$ws.Range("J141").Value = 69999.5
$ws.Range("K141").Value = 23847
$ws.Range("L141").Value = 209998.5
$ws.Range("M141").Value = -18667
$ws.Range("N141").Value = -220358.5
$ws.Range("N106").ClearContents()
$ws.Range("M125").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 1000
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("H122").Value = 2135.5625
$ws.Range("I122").Value = 2081.25
$ws.Range("K122").Value = 6243.75
$ws.Range("M122").Value = -3793.75
$ws.Range("H132").Value = 2555.2856
$ws.Range("I132").Value = 1577.6
$ws.Range("J132").Value = 4999.5
$ws.Range("K132").Value = 4732.799999999999
$ws.Range("L132").Value = 14998.5
$ws.Range("M132").Value = -2202.799999999999
$ws.Range("N132").Value = -20058.5
$ws.Range("N21").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("H20").Value = 3538.5
$ws.Range("I20").Value = 3538.5
$ws.Range("K20").Value = 3538.5
$ws.Range("M20").Value = -3291.5
$ws.Range("H86").Value = 658.61536
$ws.Range("I86").Value = 538.8
$ws.Range("J86").Value = 1058
$ws.Range("K86").Value = 538.8
$ws.Range("L86").Value = 1058
$ws.Range("M86").Value = 584.2
$ws.Range("N86").Value = -3304
$ws.Range("H89").Value = 658.61536
$ws.Range("I89").Value = 538.8
$ws.Range("J89").Value = 1058
$ws.Range("K89").Value = 2694
$ws.Range("L89").Value = 5290
$ws.Range("M89").Value = 2922
$ws.Range("N89").Value = -16522
$ws.Range("H99").Value = 2316.3333
$ws.Range("I99").Value = 2316.3333
$ws.Range("K99").Value = 2316.3333
$ws.Range("M99").Value = -818.3332999999998
$ws.Range("N19").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("H51").Value = 17866.666
$ws.Range("J51").Value = 17866.666
$ws.Range("L51").Value = 17866.666
$ws.Range("N51").Value = -19338.666
$ws.Range("H61").Value = 17866.666
$ws.Range("J61").Value = 17866.666
$ws.Range("L61").Value = 17866.666
$ws.Range("N61").Value = -18562.666
$ws.Range("H86").Value = 3500
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 3500
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -36232
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("H105").Value = 3383
$ws.Range("I105").Value = 2824.5
$ws.Range("K105").Value = 2824.5
$ws.Range("M105").Value = -1077.5
$ws.Range("N9").ClearContents()
$ws.Range("M103").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 353546.4
$ws.Range("I4").Value = 357735.06
$ws.Range("K4").Value = 1073205.18
$ws.Range("M4").Value = -1073093.18
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 1000
$ws.Range("K5").Value = 3000
$ws.Range("M5").Value = -2888
$ws.Range("H14").Value = 794.2
$ws.Range("I14").Value = 794.2
$ws.Range("K14").Value = 2382.6
$ws.Range("M14").Value = -2209.6
$ws.Range("H116").Value = 1228.6666
$ws.Range("I116").Value = 843
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 2529
$ws.Range("L116").Value = 6000
$ws.Range("M116").Value = 913
$ws.Range("N116").Value = -12884
$ws.Range("H118").Value = 364.5
$ws.Range("I118").Value = 364.5
$ws.Range("K118").Value = 1093.5
$ws.Range("M118").Value = 149.5
$ws.Range("H135").Value = 1000
$ws.Range("I135").Value = 1000
$ws.Range("K135").Value = 9000
$ws.Range("M135").Value = -6465

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 119.57143
$ws.Range("I2").Value = 110.5
$ws.Range("J2").Value = 142.25
$ws.Range("K2").Value = 110.5
$ws.Range("L2").Value = 142.25
$ws.Range("M2").Value = 2.5
$ws.Range("N2").Value = -368.25
$ws.Range("H43").Value = 9403.4
$ws.Range("I43").Value = 2017
$ws.Range("K43").Value = 2017
$ws.Range("M43").Value = -1866
$ws.Range("H57").Value = 16410.8
$ws.Range("J57").Value = 19999.75
$ws.Range("L57").Value = 19999.75
$ws.Range("N57").Value = -21639.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9166.666999999999
$ws.Range("I7").Value = 9500
$ws.Range("K7").Value = 9500
$ws.Range("M7").Value = -9388
$ws.Range("H40").Value = 3845.7693
$ws.Range("I40").Value = 3555.111
$ws.Range("J40").Value = 4499.75
$ws.Range("K40").Value = 3555.111
$ws.Range("L40").Value = 4499.75
$ws.Range("M40").Value = -3419.111
$ws.Range("N40").Value = -4771.75
$ws.Range("H46").Value = 3935.3333
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 3935.3333
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 3935.3333
$ws.Range("N46").Value = -4311.3333
$ws.Range("H61").Value = 4498
$ws.Range("I61").Value = 4498
$ws.Range("K61").Value = 4498
$ws.Range("M61").Value = -4296
$ws.Range("H69").Value = 55000
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("H72").Value = 55000
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("H113").Value = 4498
$ws.Range("I113").Value = 4498
$ws.Range("K113").Value = 4498
$ws.Range("M113").Value = -2328
$ws.Range("H126").Value = 9166.666999999999
$ws.Range("I126").Value = 9500
$ws.Range("K126").Value = 28500
$ws.Range("M126").Value = -26030
$ws.Range("M46").ClearContents()
$ws.Range("M69").ClearContents()
$ws.Range("M72").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("H126").Value = 2249
$ws.Range("I126").Value = 2998.3333
$ws.Range("K126").Value = 8994.999899999999
$ws.Range("M126").Value = -6524.999899999999
$ws.Range("H132").Value = 3833.1667
$ws.Range("I132").Value = 5499.5
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 16498.5
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -13968.5
$ws.Range("N132").Value = -14060
$ws.Range("M51").ClearContents()
